$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Question numbering column (A) changes ---
# A4 becomes a literal value (drops the shared formula that used to start there)
$ws.Range("A4").Value = 3

# A5 becomes the new start of the "fill down" formula (A = previous row + 1)
$ws.Range("A5:A6").Formula = "=A4+1"

# A9:A16 keep counting up from A8 (A7/A8 remain plain literal numbers, untouched)
$ws.Range("A9:A16").Formula = "=A8+1"

# Row 8 (A8) becomes a literal 6 instead of 7; downstream formulas recompute
# automatically through the dependency chain
$ws.Range("A8").Value = 6

# Row 17 (A17) loses its formula entirely and becomes a literal value
$ws.Range("A17").Value = 15

# Row 18 (A18) literal number also updates to 15
$ws.Range("A18").Value = 15

# --- Remove the old totals row 19 (E19 SUM formula) ---
$ws.Range("E19").ClearContents()

# --- Remove the per-model summary data in rows 21-23 (keep styled, blank C cells) ---
$ws.Range("C21:E21").ClearContents()
$ws.Range("C22:E22").ClearContents()
$ws.Range("C23:E23").ClearContents()

# --- Update the remembered selection / active cell ---
$ws.Range("G13").Select()
